$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add three new animation rows (finch, groggy, reaction)
$ws.Range("C14").Value = 10
$ws.Range("D14").Value = "finch"
$ws.Range("E14").Value = 1101
$ws.Range("F14").Value = 1131

$ws.Range("C15").Value = 11
$ws.Range("D15").Value = "groggy"
$ws.Range("E15").Value = 1132
$ws.Range("F15").Value = 1212

$ws.Range("C16").Value = 12
$ws.Range("D16").Value = "reaction"
$ws.Range("E16").Value = 1213
$ws.Range("F16").Value = 1420

# Update the active selection to match the new last entry
$ws.Range("F16").Select()
